# Update the "想去人数" (number of people interested) counts that changed
# because the site data was regenerated (gh-pages output refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value  = 13
$wsExhibit.Range("F8").Value  = 117
$wsExhibit.Range("F12").Value = 132
$wsExhibit.Range("F13").Value = 11374

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 95

# Sheet "全部类型" (all types - combined view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 95
$wsAll.Range("F9").Value  = 13
$wsAll.Range("F10").Value = 117
$wsAll.Range("F14").Value = 132
$wsAll.Range("F15").Value = 11374
